# Applies the YGA2018-conference-abstract-template.docx edit described by the diff:
#   1. Drop "respective " from "...please send it to the respective conference
#      secretaries:" so it reads "...please send it to the conference secretaries:".
#   2. Word keeps a single "_GoBack" bookmark that always marks the location of the
#      most recent edit. After the text edit above it therefore moves from the old
#      (now stale) empty paragraph to the empty paragraph right after the
#      "...Macau: LI Haifeng (haifengli@umac.mo)." line - re-adding a bookmark
#      named "_GoBack" automatically relocates (and removes the old) one, since a
#      document can only contain one "_GoBack" bookmark.

$d = $word.ActiveDocument

# --- Step 1: "the respective conference" -> "the conference" ---------------
$editRange = $d.Content
$editRange.Find.Execute("the respective ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "the ", 2)

# --- Step 2: relocate "_GoBack" to the paragraph after the e-mail list -----
$i = 0
$targetParaIndex = -1
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Macau: LI Haifeng*") {
        $targetParaIndex = $i + 1
    }
}
if ($targetParaIndex -gt 0) {
    $d.Bookmarks.Add("_GoBack", $d.Paragraphs($targetParaIndex).Range)
}
